$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "58.906.15"
$ws.Range("E2").Value = "  -0.85%  "

# Row 3
$ws.Range("D3").Value = "2.516.18"
$ws.Range("E3").Value = "  -0.05%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.997"
$ws.Range("E4").Value = "  -0.45%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "535.69"
$ws.Range("E5").Value = "  -0.54%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "136.75"
$ws.Range("E6").Value = "  -2.09%  "

# Row 7
$ws.Range("E7").Value = "  -0.20%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.568"
$ws.Range("E8").Value = "  +0.81%  "

# Row 9
$ws.Range("D9").Value = "2.514.04"
$ws.Range("E9").Value = "  -0.30%  "

# Row 10
$ws.Range("E10").Value = "  -1.33%  "

# Row 11
$ws.Range("E11").Value = "  -2.26%  "

# Row 12
$ws.Range("E12").Value = "  -1.78%  "

# Row 13
$ws.Range("E13").Value = "  -3.49%  "

# Row 14
$ws.Range("D14").Value = "2.938.04"
$ws.Range("E14").Value = "  -1.05%  "

# Row 15
$ws.Range("B15").Value = "Avalanche"
$ws.Range("C15").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "22.89"
$ws.Range("E15").Value = "  -2.40%  "

# Row 16
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "58.813.22"
$ws.Range("E16").Value = "  -0.87%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.0000139"
$ws.Range("E17").Value = "  -2.08%  "

# Row 18
$ws.Range("D18").Value = "2.512.52"
$ws.Range("E18").Value = "  -0.21%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "11.07"
$ws.Range("E19").Value = "  -0.45%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.27"
$ws.Range("E20").Value = "  -1.07%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.32"
$ws.Range("E21").Value = "  -1.04%  "

# Row 22
$ws.Range("E22").Value = "  -0.03%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.94"
$ws.Range("E23").Value = "  +1.63%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "65.68"
$ws.Range("E24").Value = "  +3.82%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.420"
$ws.Range("E25").Value = "  -0.99%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.164"
$ws.Range("E26").Value = "  -2.09%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.999"
$ws.Range("E27").Value = "  -0.15%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.51"
$ws.Range("E28").Value = "  -4.09%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.68"
$ws.Range("E29").Value = "  -3.19%  "

# Row 30
$ws.Range("D30").Value = "0.0₃0770"
$ws.Range("E30").Value = "  -1.55%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.77"
$ws.Range("E31").Value = "  -1.49%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "166.07"
$ws.Range("E32").Value = "  +0.60%  "

# Row 33
$ws.Range("E33").Value = "  +4.21%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.997"
$ws.Range("E34").Value = "  -0.14%  "

# Row 35
$ws.Range("E35").Value = "  -0.56%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "18.40"
$ws.Range("E36").Value = "  -0.67%  "

# Row 37
$ws.Range("E37").Value = "  -4.60%  "

# Row 38
$ws.Range("E38").Value = "  -3.75%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "36.70"
$ws.Range("E39").Value = "  -0.70%  "

# Row 40
$ws.Range("E40").Value = "  -0.11%  "

# Row 41
$ws.Range("E41").Value = "  -2.53%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "283.73"
$ws.Range("E42").Value = "  +1.15%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.21"
$ws.Range("E43").Value = "  -0.96%  "

# Row 44
$ws.Range("B44").Value = "Aave"
$ws.Range("C44").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "132.45"
$ws.Range("E44").Value = "  +6.40%  "

# Row 45
$ws.Range("B45").Value = "FirstDigitalUSD"
$ws.Range("C45").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.999"
$ws.Range("E45").Value = "  +0.10%  "

# Row 46
$ws.Range("E46").Value = "  +0.74%  "

# Row 47
$ws.Range("E47").Value = "  -0.15%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.0924"
$ws.Range("E48").Value = "  -1.19%  "

# Row 49
$ws.Range("E49").Value = "  -2.27%  "

# Row 50
$ws.Range("E50").Value = "  -2.58%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.22"
$ws.Range("E51").Value = "  -3.76%  "
